# Sync tweaks to the POM package:
#  - drop the stale "devices" sheet and promote "devices (2)" to "devices"
#  - refresh the device-matrix values on the (now canonical) "devices" sheet
#  - leave "signIn" / "findPost" sheets untouched

$wb = $excel.ActiveWorkbook

# Remove the old "devices" sheet (its data is superseded by "devices (2)").
$wb.Worksheets.Item("devices").Delete() | Out-Null

# Promote "devices (2)" to be the canonical "devices" sheet.
$ws = $wb.Worksheets.Item("devices (2)")
$ws.Name = "devices"

# Update the device-matrix cell values.
$ws.Range("I2").Value = "Avner|Raj|Brian"

$ws.Range("E3").ClearContents() | Out-Null
$ws.Range("D3").Value = "Samsung"
$ws.Range("I3").Value = "Avner|raj|Brian"

$ws.Range("A4").Value = "chrome"
$ws.Range("A5").Value = "firefox"

$ws.Range("D6").ClearContents() | Out-Null
$ws.Range("E6").Value = "iPad*.*"
$ws.Range("I6").Value = "Shared"

# Make "devices" the active/selected sheet and park the selection at J7.
$ws.Activate()
$ws.Range("J7").Select() | Out-Null
